# Weekly update: insert 3 new "Berenjena" (Vega Central Mapocho de Santiago)
# price rows for the newest reporting date, pushing the previously-existing
# rows (214..256) down to (217..259).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows above the current row 214, shifting everything below
# (including the previous row 214) down by 3 rows.
$ws.Rows.Item(214).Resize(3).Insert()

# ---- New row 214 ----
$ws.Range("A214").Value = 9
$ws.Range("B214").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C214").Value = "Metropolitana"
$ws.Range("D214").Value = 44644
$ws.Range("E214").Value = 13
$ws.Range("F214").Value = 100112001
$ws.Range("G214").Value = "Berenjena"
$ws.Range("H214").Value = "Sin especificar"
$ws.Range("I214").Value = "Primera"
$ws.Range("J214").Value = 97
$ws.Range("K214").Value = 8000
$ws.Range("L214").Value = 9000
$ws.Range("M214").Value = 8495
$ws.Range("N214").Value = "$/caja 50 unidades"
$ws.Range("O214").Value = "Región de Arica y Parinacota"
$ws.Range("P214").Value = 170
$ws.Range("Q214").Value = 50
$ws.Range("R214").Value = "Hortaliza"

# ---- New row 215 ----
$ws.Range("A215").Value = 9
$ws.Range("B215").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C215").Value = "Metropolitana"
$ws.Range("D215").Value = 44644
$ws.Range("E215").Value = 13
$ws.Range("F215").Value = 100112001
$ws.Range("G215").Value = "Berenjena"
$ws.Range("H215").Value = "Sin especificar"
$ws.Range("I215").Value = "Primera"
$ws.Range("J215").Value = 80
$ws.Range("K215").Value = 10000
$ws.Range("L215").Value = 10000
$ws.Range("M215").Value = 10000
$ws.Range("N215").Value = "$/caja 60 unidades"
$ws.Range("O215").Value = "Región de O'Higgins"
$ws.Range("P215").Value = 167
$ws.Range("Q215").Value = 60
$ws.Range("R215").Value = "Hortaliza"

# ---- New row 216 ----
$ws.Range("A216").Value = 9
$ws.Range("B216").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C216").Value = "Metropolitana"
$ws.Range("D216").Value = 44644
$ws.Range("E216").Value = 13
$ws.Range("F216").Value = 100112001
$ws.Range("G216").Value = "Berenjena"
$ws.Range("H216").Value = "Sin especificar"
$ws.Range("I216").Value = "Primera"
$ws.Range("J216").Value = 20
$ws.Range("K216").Value = 12000
$ws.Range("L216").Value = 12000
$ws.Range("M216").Value = 12000
$ws.Range("N216").Value = "$/caja 70 unidades"
$ws.Range("O216").Value = "Limache"
$ws.Range("P216").Value = 171
$ws.Range("Q216").Value = 70
$ws.Range("R216").Value = "Hortaliza"
